# Update the "grupos" (groups) standings table on the active sheet.
# Groups B and C had extra match results recorded, changing Po (wins),
# Jo/GP (goals for), GC (goals against), Saldo (goal diff) and Pontos
# (points) for several teams.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Group B rows ---------------------------------------------------

# Row 6: Vasquinho
$ws.Range("C6").Value = 3
$ws.Range("D6").Value = 10
$ws.Range("E6").Value = 13
$ws.Range("F6").Value = -3
$ws.Range("G6").Value = 3

# Row 7: Es. Po. Seguro
$ws.Range("C7").Value = 9
$ws.Range("D7").Value = 13
$ws.Range("E7").Value = 6
$ws.Range("F7").Value = 7
$ws.Range("G7").Value = 3

# Row 8: Lagoa Verde
$ws.Range("D8").Value = 10
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 4
$ws.Range("G8").Value = 3

# Row 9: Pró Limp
$ws.Range("D9").Value = 8
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = -7
$ws.Range("G9").Value = 3

# --- Group C rows ------------------------------------------------------

# Row 11: At. Ma. Martins
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = 7
$ws.Range("E11").Value = 8
$ws.Range("F11").Value = -1
$ws.Range("G11").Value = 2

# Row 12: MEC Futsal
$ws.Range("D12").Value = 5
$ws.Range("E12").Value = 12
$ws.Range("F12").Value = -7
$ws.Range("G12").Value = 2

# --- Column E width now differs from column D (it got its own, slightly
# wider, best-fit column) instead of sharing a combined D:E column band.
# ColumnWidth uses the "characters" unit, which is offset from the raw
# saved column "width" (XML) by the standard 5/6 char padding; 3.1666...
# (= 4 - 5/6) round-trips to a saved width of exactly 4.
$ws.Columns.Item(5).ColumnWidth = 3.1666666666666665

# --- View state: the sheet was re-selected/zoomed before saving ---------
$ws.Range("J14").Select()
$excel.ActiveWindow.Zoom = 110
